# Add a "name" column at the front of the sheet (logging needs a way to
# identify which wallet/account a row belongs to), plus small layout fixes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; existing columns (address, private,
# proxy, okx_api, ref_link) shift right from A-E to B-F.
$ws.Columns.Item(1).Insert()

# --- values -----------------------------------------------------------
$ws.Range("A1").Value = "name"
$ws.Range("A2").Value = "name1"
$ws.Range("A3").Value = 0

# --- header formatting --------------------------------------------------
# A1 should look like the rest of the header row (bold 14pt, bordered,
# center/top aligned) - just clone that existing format onto it.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- "name" column data formatting --------------------------------------
# Center the name values (A2:A3) both horizontally and vertically, keeping
# the default font.
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4108
$ws.Range("A2").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- other data columns formatting (B:F, rows 2-3) ----------------------
# Same center/center alignment, keeping the existing 12pt font.
$ws.Range("B2").HorizontalAlignment = -4108
$ws.Range("B2").VerticalAlignment = -4108
$ws.Range("B2").Copy()
$ws.Range("B2:F3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- column widths --------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 21.6666666666667
$ws.Columns.Item(2).ColumnWidth = 16.1666666666667
$ws.Columns.Item(3).ColumnWidth = 14.3333333333333
$ws.Columns.Item(4).ColumnWidth = 22
$ws.Columns.Item(5).ColumnWidth = 36
$ws.Columns.Item(6).ColumnWidth = 39.6666666666667

# --- selection ------------------------------------------------------------
$ws.Range("D12").Select()
